$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D retains its original text formatting, since many price
# values look numeric (e.g. "15.46") and Excel would otherwise silently
# coerce them into real numbers, dropping formatting like trailing zeros.

$ws.Range('D2').Value = '26.358.67'
$ws.Range('E2').Value = '  -0.61%  '

$ws.Range('D3').Value = '1.715.62'
$ws.Range('E3').Value = '  -1.21%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9974'
$ws.Range('E4').Value = '  -0.21%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.93'
$ws.Range('E5').Value = '  -2.41%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9977'
$ws.Range('E6').Value = '  -0.20%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4865'
$ws.Range('E7').Value = '  -0.60%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2584'
$ws.Range('E8').Value = '  -3.06%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06163'
$ws.Range('E9').Value = '  -3.25%  '

$ws.Range('D10').Value = '1.722.04'
$ws.Range('E10').Value = '  -0.75%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.06958'
$ws.Range('E11').Value = '  -1.15%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.46'
$ws.Range('E12').Value = '  -1.53%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.486'
$ws.Range('E13').Value = '  -2.48%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5965'
$ws.Range('E14').Value = '  -2.22%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '76.45'
$ws.Range('E15').Value = '  -1.26%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9975'
$ws.Range('E16').Value = '  -0.22%  '

$ws.Range('D17').Value = '26.355.45'
$ws.Range('E17').Value = '  -0.58%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9976'
$ws.Range('E18').Value = '  -0.22%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007087'
$ws.Range('E19').Value = '  -4.52%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.23'
$ws.Range('E20').Value = '  -2.48%  '

$ws.Range('D21').Value = '1.943.00'
$ws.Range('E21').Value = '  -0.50%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.404'
$ws.Range('E22').Value = '  -3.76%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.442'
$ws.Range('E23').Value = '  -3.04%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.049'
$ws.Range('E24').Value = '  -3.55%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '136.58'
$ws.Range('E25').Value = '  -3.00%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.20'
$ws.Range('E26').Value = '  -1.60%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.412'
$ws.Range('E27').Value = '  -0.16%  '

$ws.Range('B28').Value = 'BitcoinCash'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '105.77'
$ws.Range('E28').Value = '  -2.09%  '

$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.729'
$ws.Range('E29').Value = '  -2.13%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.872'
$ws.Range('E30').Value = '  -3.98%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.07952'
$ws.Range('E31').Value = '  -1.06%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.607'
$ws.Range('E32').Value = '  -2.98%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04426'
$ws.Range('E33').Value = '  -3.44%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.610'
$ws.Range('E34').Value = '  +0.00%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9942'
$ws.Range('E35').Value = '  -1.48%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6172'
$ws.Range('E36').Value = '  -3.03%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9332'
$ws.Range('E37').Value = '  +4.22%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.972'
$ws.Range('E38').Value = '  -2.08%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.373'
$ws.Range('E39').Value = '  -1.14%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9966'
$ws.Range('E40').Value = '  -0.76%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.01471'
$ws.Range('E41').Value = '  -2.15%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '99.16'
$ws.Range('E42').Value = '  -3.84%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.436'
$ws.Range('E43').Value = '  +0.67%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3803'
$ws.Range('E44').Value = '  -2.24%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '6.820'
$ws.Range('E45').Value = '  -0.98%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1150'
$ws.Range('E46').Value = '  -3.01%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05350'
$ws.Range('E47').Value = '  -0.78%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '30.37'
$ws.Range('E48').Value = '  -0.46%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.736'
$ws.Range('E49').Value = '  -0.74%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '51.13'
$ws.Range('E50').Value = '  -1.18%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.214'
$ws.Range('E51').Value = '  -3.78%  '
